# "updated activity till excel form"
# Robin Uthappa (Rajasthan Royals) match-by-match batting log:
#   - existing rows 2-12 get their runs/balls/fours/sixes reshuffled
#     (latest match data refreshed from source), and
#   - a new match (row 13) is appended.
#
# Columns C (runs), D (balls), E (fours) and F (sixes) hold numeric-looking
# values that must stay stored as TEXT (matching the workbook's existing
# "number stored as text" convention), so they are written with a leading
# apostrophe to force Excel to keep them as strings instead of coercing to
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $ws.Range($range).Value = "'" + $value
}

# Only cells whose text actually changes are touched below - cells whose
# value is unchanged between before/after (e.g. E2, F3, F4, F6, E7, F7, F9,
# F10, F11, F12) are intentionally left alone.

# Row 2
Set-TextValue "C2" "30"
Set-TextValue "D2" "23"
Set-TextValue "F2" "2"

# Row 3
Set-TextValue "C3" "19"
Set-TextValue "D3" "13"
Set-TextValue "E3" "2"

# Row 4
Set-TextValue "C4" "13"
Set-TextValue "D4" "11"
Set-TextValue "E4" "2"

# Row 5
Set-TextValue "C5" "5"
Set-TextValue "D5" "9"
Set-TextValue "E5" "0"
Set-TextValue "F5" "0"

# Row 6
Set-TextValue "C6" "6"
Set-TextValue "D6" "2"
Set-TextValue "E6" "0"

# Row 7
Set-TextValue "C7" "4"
Set-TextValue "D7" "9"

# Row 8
Set-TextValue "C8" "32"
Set-TextValue "D8" "27"
Set-TextValue "E8" "3"
Set-TextValue "F8" "1"

# Row 9
Set-TextValue "C9" "18"
Set-TextValue "D9" "15"
Set-TextValue "E9" "1"

# Row 10
Set-TextValue "C10" "17"
Set-TextValue "D10" "22"
Set-TextValue "E10" "1"

# Row 11
Set-TextValue "C11" "41"
Set-TextValue "D11" "22"
Set-TextValue "E11" "7"

# Row 12
Set-TextValue "C12" "9"
Set-TextValue "D12" "4"
Set-TextValue "E12" "2"

# Row 13 (new match row appended at the bottom)
# Note: the player name cells elsewhere in this sheet (and the sheet/file
# name itself) end with a non-breaking space (U+00A0), not a plain space -
# match that here so the new row's string reuses the same text run.
$ws.Range("A13").Value = "Robin Uthappa "
$ws.Range("B13").Value = "Rajasthan Royals"
Set-TextValue "C13" "2"
Set-TextValue "D13" "7"
Set-TextValue "E13" "0"
Set-TextValue "F13" "0"
